$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove two rows that are no longer needed so the remaining detail
# --- rows (and the footer block below them) shift up to their target
# --- row numbers (79 -> 77, 84/85 -> 82/83), while the special
# --- bottom-border style that lived on row 79 rides along to the new
# --- last data row automatically.
$ws.Rows.Item(78).Delete()
$ws.Rows.Item(16).Delete()

# --- Rewrite the detail table (rows 16-77): Tipo Doc, N Doc, Nombre,
# --- Periodo Mora, Valor Mora, Salario Basico
$data = @(
    @(16, "CC", "1128047433", "FRANKLIN BARRAGAN ECHEVERRIA", "1610", 5344, 1336000),
    @(17, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1808", 31249, 1610185),
    @(18, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1809", 31249, 1610185),
    @(19, "CC", "1076820808", "WALTER AGUILAR MOSQUERA", "1809", 44213, 1105310),
    @(20, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1810", 31249, 1610185),
    @(21, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1811", 31249, 1610185),
    @(22, "CC", "1076820808", "WALTER AGUILAR MOSQUERA", "1811", 44213, 1105310),
    @(23, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1812", 31249, 1610185),
    @(24, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1901", 31249, 1610185),
    @(25, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1902", 31249, 1610185),
    @(26, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1903", 31249, 1610185),
    @(27, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1904", 31249, 1610185),
    @(28, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1905", 31249, 1610185),
    @(29, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1906", 31249, 1610185),
    @(30, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1907", 31249, 1610185),
    @(31, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1908", 31249, 1610185),
    @(32, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1909", 31249, 1610185),
    @(33, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1910", 31249, 1610185),
    @(34, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1911", 31249, 1610185),
    @(35, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "1912", 31249, 1610185),
    @(36, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2001", 31249, 1610185),
    @(37, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2002", 31249, 1610185),
    @(38, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2003", 31249, 1610185),
    @(39, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2004", 31249, 1610185),
    @(40, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2005", 31249, 1610185),
    @(41, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2006", 31249, 1610185),
    @(42, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2007", 31249, 1610185),
    @(43, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2008", 31249, 1610185),
    @(44, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2009", 31249, 1610185),
    @(45, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2010", 31249, 1610185),
    @(46, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2011", 31249, 1610185),
    @(47, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2012", 31249, 1610185),
    @(48, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2101", 31249, 1610185),
    @(49, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2102", 31249, 1610185),
    @(50, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2103", 31249, 1610185),
    @(51, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2104", 31249, 1610185),
    @(52, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2105", 31249, 1610185),
    @(53, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2106", 31249, 1610185),
    @(54, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2107", 31249, 1610185),
    @(55, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2108", 31249, 1610185),
    @(56, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2109", 31249, 1610185),
    @(57, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2110", 31249, 1610185),
    @(58, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2111", 31249, 1610185),
    @(59, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2112", 31249, 1610185),
    @(60, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2201", 31249, 1610185),
    @(61, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2202", 31249, 1610185),
    @(62, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2203", 40000, 1610185),
    @(63, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2204", 40000, 1610185),
    @(64, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2205", 40000, 1610185),
    @(65, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2206", 40000, 1610185),
    @(66, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2207", 40000, 1610185),
    @(67, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2208", 40000, 1610185),
    @(68, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2209", 40000, 1610185),
    @(69, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2210", 40000, 1610185),
    @(70, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2211", 40000, 1610185),
    @(71, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2212", 40000, 1610185),
    @(72, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2301", 40000, 1610185),
    @(73, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2302", 40000, 1610185),
    @(74, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2303", 40000, 1610185),
    @(75, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2304", 40000, 1610185),
    @(76, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2305", 40000, 1610185),
    @(77, "CC", "7938782", "ANDRES GUILLERMO PAJARO MORALES", "2306", 31249, 1610185)

)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# --- Summary header fields
$ws.Range("E11").Value = 2068726
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 60

# --- Column D ("Nombre Trabajador") re-bestfits narrower now that the
# --- longest remaining name is shorter than before.
$ws.Columns.Item(4).AutoFit()
